$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Cells.Item(6, 8).Value = 151.16667
$ws.Cells.Item(6, 9).Value = 70.333336
$ws.Cells.Item(6, 10).Value = 232
$ws.Cells.Item(6, 11).Value = 211.000008
$ws.Cells.Item(6, 12).Value = 696
$ws.Cells.Item(6, 13).Value = -99.00000800000001
$ws.Cells.Item(6, 14).Value = -920
# Row 9
$ws.Cells.Item(9, 8).Value = 1182.4375
$ws.Cells.Item(9, 9).Value = 1518.25
$ws.Cells.Item(9, 10).Value = 175
$ws.Cells.Item(9, 11).Value = 1518.25
$ws.Cells.Item(9, 12).Value = 175
$ws.Cells.Item(9, 13).Value = -1349.25
$ws.Cells.Item(9, 14).Value = -513
# Row 31
$ws.Cells.Item(31, 8).Value = 1127.4
$ws.Cells.Item(31, 9).Value = 1127.4
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 3382.2
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -3152.2
# Row 62
$ws.Cells.Item(62, 8).Value = 5000
$ws.Cells.Item(62, 9).Value = 5000
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 5000
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -4376
# Row 65
$ws.Cells.Item(65, 8).Value = 5000
$ws.Cells.Item(65, 9).Value = 5000
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 25000
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -21880
# Row 106
$ws.Cells.Item(106, 8).Value = 38365
$ws.Cells.Item(106, 9).Value = 38365
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 11).Value = 38365
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = -37734
# Row 111
$ws.Cells.Item(111, 8).Value = 49666.5
$ws.Cells.Item(111, 9).Value = 49500
$ws.Cells.Item(111, 10).Value = 49999.5
$ws.Cells.Item(111, 11).Value = 148500
$ws.Cells.Item(111, 12).Value = 149998.5
$ws.Cells.Item(111, 13).Value = -145433
$ws.Cells.Item(111, 14).Value = -156132.5
# Row 131
$ws.Cells.Item(131, 8).Value = 14718
$ws.Cells.Item(131, 9).Value = 15091.286
$ws.Cells.Item(131, 10).Value = 12105
$ws.Cells.Item(131, 11).Value = 45273.858
$ws.Cells.Item(131, 12).Value = 36315
$ws.Cells.Item(131, 13).Value = -40233.858
$ws.Cells.Item(131, 14).Value = -46395
# Row 137
$ws.Cells.Item(137, 8).Value = 1430.6364
$ws.Cells.Item(137, 9).Value = 1304.1111
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 3912.3333
$ws.Cells.Item(137, 12).Value = 6000
$ws.Cells.Item(137, 13).Value = -1362.3333
$ws.Cells.Item(137, 14).Value = -11100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2335612.5
$ws.Cells.Item(32, 9).Value = 2189850.5
$ws.Cells.Item(32, 10).Value = 7000000
$ws.Cells.Item(32, 11).Value = 2189850.5
$ws.Cells.Item(32, 12).Value = 7000000
$ws.Cells.Item(32, 13).Value = -2189563.5
$ws.Cells.Item(32, 14).Value = -7000574
# Row 61
$ws.Cells.Item(61, 8).Value = 2590.3635
$ws.Cells.Item(61, 9).Value = 2349.5
$ws.Cells.Item(61, 10).Value = 4999
$ws.Cells.Item(61, 11).Value = 2349.5
$ws.Cells.Item(61, 12).Value = 4999
$ws.Cells.Item(61, 13).Value = -2137.5
$ws.Cells.Item(61, 14).Value = -5423
# Row 74
$ws.Cells.Item(74, 8).Value = 2710.5
$ws.Cells.Item(74, 9).Value = 3518.3333
$ws.Cells.Item(74, 10).Value = 287
$ws.Cells.Item(74, 11).Value = 3518.3333
$ws.Cells.Item(74, 12).Value = 287
$ws.Cells.Item(74, 13).Value = -2644.3333
$ws.Cells.Item(74, 14).Value = -2035
# Row 77
$ws.Cells.Item(77, 8).Value = 2710.5
$ws.Cells.Item(77, 9).Value = 3518.3333
$ws.Cells.Item(77, 10).Value = 287
$ws.Cells.Item(77, 11).Value = 17591.6665
$ws.Cells.Item(77, 12).Value = 1435
$ws.Cells.Item(77, 13).Value = -13223.6665
$ws.Cells.Item(77, 14).Value = -10171
# Row 102
$ws.Cells.Item(102, 8).Value = 1302.75
$ws.Cells.Item(102, 9).Value = 1302.75
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1302.75
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = 319.25
# Row 104
$ws.Cells.Item(104, 8).Value = 29000
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 29000
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 29000
$ws.Cells.Item(104, 14).Value = -35988
# Row 136
$ws.Cells.Item(136, 8).Value = 2590.3635
$ws.Cells.Item(136, 9).Value = 2349.5
$ws.Cells.Item(136, 10).Value = 4999
$ws.Cells.Item(136, 11).Value = 7048.5
$ws.Cells.Item(136, 12).Value = 14997
$ws.Cells.Item(136, 13).Value = -4498.5
$ws.Cells.Item(136, 14).Value = -20097

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 402.1
$ws.Cells.Item(94, 9).Value = 391.22223
$ws.Cells.Item(94, 10).Value = 500
$ws.Cells.Item(94, 11).Value = 391.22223
$ws.Cells.Item(94, 12).Value = 500
$ws.Cells.Item(94, 13).Value = 59.77776999999998
$ws.Cells.Item(94, 14).Value = -1402
# Row 105
$ws.Cells.Item(105, 8).Value = 2916.6667
$ws.Cells.Item(105, 9).Value = 2625
$ws.Cells.Item(105, 10).Value = 3500
$ws.Cells.Item(105, 11).Value = 2625
$ws.Cells.Item(105, 12).Value = 3500
$ws.Cells.Item(105, 13).Value = -878
$ws.Cells.Item(105, 14).Value = -6994
# Row 107
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).ClearContents()
# Row 11
$ws.Cells.Item(11, 8).Value = 1792.75
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 1792.75
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(11, 12).Value = 1792.75
$ws.Cells.Item(11, 14).Value = -2072.75
$ws.Cells.Item(11, 13).ClearContents()
# Row 69
$ws.Cells.Item(69, 8).Value = 40000
$ws.Cells.Item(69, 9).Value = 40000
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 40000
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = -39251
# Row 72
$ws.Cells.Item(72, 8).Value = 40000
$ws.Cells.Item(72, 9).Value = 40000
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 120000
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = -116256
# Row 99
$ws.Cells.Item(99, 8).Value = 1492.48
$ws.Cells.Item(99, 9).Value = 1251
$ws.Cells.Item(99, 10).Value = 1854.7
$ws.Cells.Item(99, 11).Value = 1251
$ws.Cells.Item(99, 12).Value = 1854.7
$ws.Cells.Item(99, 13).Value = 247
$ws.Cells.Item(99, 14).Value = -4850.7
# Row 107
$ws.Cells.Item(107, 9).Value = 1499
$ws.Cells.Item(107, 10).Value = 1806.5
$ws.Cells.Item(107, 11).Value = 1499
$ws.Cells.Item(107, 12).Value = 1806.5
$ws.Cells.Item(107, 13).Value = 421
$ws.Cells.Item(107, 14).Value = -5646.5
# Row 126
$ws.Cells.Item(126, 8).Value = 1492.48
$ws.Cells.Item(126, 9).Value = 1251
$ws.Cells.Item(126, 10).Value = 1854.7
$ws.Cells.Item(126, 11).Value = 3753
$ws.Cells.Item(126, 12).Value = 5564.1
$ws.Cells.Item(126, 13).Value = -1283
$ws.Cells.Item(126, 14).Value = -10504.1

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 31
$ws.Cells.Item(2, 9).Value = 18.75
$ws.Cells.Item(2, 10).Value = 47.333332
$ws.Cells.Item(2, 11).Value = 112.5
$ws.Cells.Item(2, 12).Value = 283.999992
$ws.Cells.Item(2, 13).Value = 0.5
$ws.Cells.Item(2, 14).Value = -509.999992
# Row 34
$ws.Cells.Item(34, 8).Value = 1250
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 1250
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 3750
$ws.Cells.Item(34, 14).Value = -3918
$ws.Cells.Item(34, 13).ClearContents()
# Row 38
$ws.Cells.Item(38, 8).Value = 78
$ws.Cells.Item(38, 9).Value = 77.666664
$ws.Cells.Item(38, 10).Value = 79.5
$ws.Cells.Item(38, 11).Value = 232.999992
$ws.Cells.Item(38, 12).Value = 238.5
$ws.Cells.Item(38, 13).Value = 114.000008
$ws.Cells.Item(38, 14).Value = -932.5
# Row 113
$ws.Cells.Item(113, 8).Value = 1252.9
$ws.Cells.Item(113, 9).Value = 1498.3334
$ws.Cells.Item(113, 10).Value = 1147.7142
$ws.Cells.Item(113, 11).Value = 4495.0002
$ws.Cells.Item(113, 12).Value = 3443.1426
$ws.Cells.Item(113, 13).Value = -2325.0002
$ws.Cells.Item(113, 14).Value = -7783.142599999999
# Row 121
$ws.Cells.Item(121, 8).Value = 2581.6
$ws.Cells.Item(121, 9).Value = 628.3333
$ws.Cells.Item(121, 10).Value = 3418.7144
$ws.Cells.Item(121, 11).Value = 1884.9999
$ws.Cells.Item(121, 12).Value = 10256.1432
$ws.Cells.Item(121, 13).Value = -574.9999
$ws.Cells.Item(121, 14).Value = -12876.1432
# Row 131
$ws.Cells.Item(131, 8).Value = 2440.1052
$ws.Cells.Item(131, 9).Value = 1483
$ws.Cells.Item(131, 10).Value = 2998.4167
$ws.Cells.Item(131, 11).Value = 4449
$ws.Cells.Item(131, 12).Value = 8995.250100000001
$ws.Cells.Item(131, 13).Value = 591
$ws.Cells.Item(131, 14).Value = -19075.2501
# Row 139
$ws.Cells.Item(139, 8).Value = 2097.3333
$ws.Cells.Item(139, 9).Value = 1146.5
$ws.Cells.Item(139, 10).Value = 3999
$ws.Cells.Item(139, 11).Value = 3439.5
$ws.Cells.Item(139, 12).Value = 11997
$ws.Cells.Item(139, 13).Value = 1700.5
$ws.Cells.Item(139, 14).Value = -22277

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Cells.Item(11, 8).Value = 4607200
$ws.Cells.Item(11, 9).Value = 7562500
$ws.Cells.Item(11, 10).Value = 666800
$ws.Cells.Item(11, 11).Value = 7562500
$ws.Cells.Item(11, 12).Value = 666800
$ws.Cells.Item(11, 13).Value = -7562361
$ws.Cells.Item(11, 14).Value = -667078
# Row 25
$ws.Cells.Item(25, 8).Value = 100000
$ws.Cells.Item(25, 9).Value = 50000
$ws.Cells.Item(25, 10).Value = 200000
$ws.Cells.Item(25, 11).Value = 50000
$ws.Cells.Item(25, 12).Value = 200000
$ws.Cells.Item(25, 13).Value = -49471
$ws.Cells.Item(25, 14).Value = -201058
# Row 113
$ws.Cells.Item(113, 8).Value = 575.0909
$ws.Cells.Item(113, 9).Value = 485
$ws.Cells.Item(113, 10).Value = 815.3333
$ws.Cells.Item(113, 11).Value = 485
$ws.Cells.Item(113, 12).Value = 815.3333
$ws.Cells.Item(113, 13).Value = 1685
$ws.Cells.Item(113, 14).Value = -5155.3333
# Row 132
$ws.Cells.Item(132, 8).Value = 6935.3335
$ws.Cells.Item(132, 9).Value = 7453
$ws.Cells.Item(132, 10).Value = 2794
$ws.Cells.Item(132, 11).Value = 22359
$ws.Cells.Item(132, 12).Value = 8382
$ws.Cells.Item(132, 13).Value = -19829
$ws.Cells.Item(132, 14).Value = -13442
# Row 136
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Cells.Item(2, 8).Value = 3400000
$ws.Cells.Item(2, 9).Value = 100000
$ws.Cells.Item(2, 10).Value = 10000000
$ws.Cells.Item(2, 11).Value = 100000
$ws.Cells.Item(2, 12).Value = 10000000
$ws.Cells.Item(2, 13).Value = -99888
$ws.Cells.Item(2, 14).Value = -10000224
# Row 23
$ws.Cells.Item(23, 8).Value = 10500
$ws.Cells.Item(23, 9).Value = 10500
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 10500
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = -10270
# Row 122
$ws.Cells.Item(122, 8).Value = 6009.963
$ws.Cells.Item(122, 9).Value = 4790.6875
$ws.Cells.Item(122, 10).Value = 7783.4546
$ws.Cells.Item(122, 11).Value = 14372.0625
$ws.Cells.Item(122, 12).Value = 23350.3638
$ws.Cells.Item(122, 13).Value = -11922.0625
$ws.Cells.Item(122, 14).Value = -28250.3638
# Row 140
$ws.Cells.Item(140, 8).Value = 29933.334
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 29933.334
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 29933.334
$ws.Cells.Item(140, 14).Value = -40293.334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Cells.Item(6, 8).Value = 3076.8
$ws.Cells.Item(6, 9).Value = 2673.5
$ws.Cells.Item(6, 10).Value = 3345.6667
$ws.Cells.Item(6, 11).Value = 2673.5
$ws.Cells.Item(6, 12).Value = 3345.6667
$ws.Cells.Item(6, 13).Value = -2558.5
$ws.Cells.Item(6, 14).Value = -3575.6667
# Row 7
$ws.Cells.Item(7, 8).Value = 17747.5
$ws.Cells.Item(7, 9).Value = 495
$ws.Cells.Item(7, 10).Value = 35000
$ws.Cells.Item(7, 11).Value = 495
$ws.Cells.Item(7, 12).Value = 35000
$ws.Cells.Item(7, 13).Value = -382
$ws.Cells.Item(7, 14).Value = -35226
# Row 12
$ws.Cells.Item(12, 8).Value = 1975
$ws.Cells.Item(12, 9).Value = 500
$ws.Cells.Item(12, 10).Value = 2466.6667
$ws.Cells.Item(12, 11).Value = 500
$ws.Cells.Item(12, 12).Value = 2466.6667
$ws.Cells.Item(12, 13).Value = -358
$ws.Cells.Item(12, 14).Value = -2750.6667
# Row 74
$ws.Cells.Item(74, 8).Value = 21903.572
$ws.Cells.Item(74, 9).Value = 20777
$ws.Cells.Item(74, 10).Value = 22354.2
$ws.Cells.Item(74, 11).Value = 20777
$ws.Cells.Item(74, 12).Value = 22354.2
$ws.Cells.Item(74, 13).Value = -19841
$ws.Cells.Item(74, 14).Value = -24226.2
# Row 77
$ws.Cells.Item(77, 8).Value = 21903.572
$ws.Cells.Item(77, 9).Value = 20777
$ws.Cells.Item(77, 10).Value = 22354.2
$ws.Cells.Item(77, 11).Value = 62331
$ws.Cells.Item(77, 12).Value = 67062.60000000001
$ws.Cells.Item(77, 13).Value = -57651
$ws.Cells.Item(77, 14).Value = -76422.60000000001
# Row 113
$ws.Cells.Item(113, 8).Value = 758.1667
$ws.Cells.Item(113, 9).Value = 749
$ws.Cells.Item(113, 10).Value = 776.5
$ws.Cells.Item(113, 11).Value = 2247
$ws.Cells.Item(113, 12).Value = 2329.5
$ws.Cells.Item(113, 13).Value = -77
$ws.Cells.Item(113, 14).Value = -6669.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1500
$ws.Cells.Item(122, 9).Value = 1500
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(122, 14).ClearContents()
# Row 126
$ws.Cells.Item(126, 8).Value = 3113.4285
$ws.Cells.Item(126, 9).Value = 1774.25
$ws.Cells.Item(126, 10).Value = 4899
$ws.Cells.Item(126, 11).Value = 5322.75
$ws.Cells.Item(126, 12).Value = 14697
$ws.Cells.Item(126, 13).Value = -2852.75
$ws.Cells.Item(126, 14).Value = -19637
# Row 136
$ws.Cells.Item(136, 8).Value = 2130.9285
$ws.Cells.Item(136, 9).Value = 2130.9285
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 6392.7855
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -3842.7855
# Row 138
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()
# Row 139
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
